# Insert a new data row at row 206 (pushing the existing row 206..275 down
# to 207..276) and populate it with a new "Betarraga" price observation for
# the "Macroferia Regional de Talca" market.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("206:206").Insert()

$ws.Range("A206").Value = 5
$ws.Range("B206").Value = "Macroferia Regional de Talca"
$ws.Range("C206").Value = "Maule"
$ws.Range("D206").Value = 44627
$ws.Range("E206").Value = 7
$ws.Range("F206").Value = 100114014
$ws.Range("G206").Value = "Betarraga"
$ws.Range("H206").Value = "Sin especificar"
$ws.Range("I206").Value = "Primera"
$ws.Range("J206").Value = 3000
$ws.Range("K206").Value = 800
$ws.Range("L206").Value = 800
$ws.Range("M206").Value = 800
$ws.Range("N206").Value = "`$/paquete 5 unidades"
$ws.Range("O206").Value = "Región del Maule"
$ws.Range("P206").Value = 160
$ws.Range("Q206").Value = 5
$ws.Range("R206").Value = "Hortaliza"
